# MirrorMe - modifications made to accommodate for the refactoring changes
#
# The "requires" relation that used to be crammed into a single comma-
# separated cell (column G: "t03, t04" / "t05, t06") is split out so each
# referenced template-text id gets its own column (G keeps the first id,
# a new column H holds the second id). This is done on both worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("MirrorMe Example Argument", "MirroMe Voorbeeldredenering")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Copy formatting from column G into the new column H for the rows
    #     that carry data, so H matches the look of its row (header rows
    #     use the "Neutral" style, data rows use the plain column style).
    $ws.Range("G1").Copy()
    $ws.Range("H1").PasteSpecial(-4122)

    $ws.Range("G2").Copy()
    $ws.Range("H2").PasteSpecial(-4122)

    $ws.Range("G5").Copy()
    $ws.Range("H5").PasteSpecial(-4122)

    $ws.Range("G7").Copy()
    $ws.Range("H7").PasteSpecial(-4122)

    $wb.Application.CutCopyMode = $false

    # --- Fix up the "[Template,]" placeholder header back to "Template".
    $ws.Range("G2").Value = "Template"

    # --- Header rows: duplicate the "requires"/"Template" header labels
    #     into the new H column.
    $ws.Range("H1").Value = "requires"
    $ws.Range("H2").Value = "Template"

    # --- Split the two "requires" list cells into one id per column.
    if ($sheetName -eq "MirrorMe Example Argument") {
        $ws.Range("G5").Value = "t03"
        $ws.Range("H5").Value = "t04"

        $ws.Range("G7").Value = "t05"
        $ws.Range("H7").Value = "t06"
    } else {
        $ws.Range("G5").Value = "TText_03"
        $ws.Range("H5").Value = "TText_04"

        $ws.Range("G7").Value = "TText_05"
        $ws.Range("H7").Value = "TText_06"
    }

    # --- Move the active selection to H7, matching where the author's
    #     cursor ended up after making the edit.
    $ws.Range("H7").Select()
}

$wb.Worksheets.Item("MirroMe Voorbeeldredenering").Activate()
